$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the columns C:E entirely (no longer used)
$ws.Range("C1:E4").Clear()

# Clear A1 header cell (the "Gender" label is removed, only B1 keeps a header now)
$ws.Range("A1").Clear()

# Update header for column B
$ws.Range("B1").Value = "Percentage of Players"

# Force column B percentage cells to be stored as plain text (not a numeric
# percentage) - set to Text format first so Excel doesn't auto-convert the
# "83.59%"-style strings into numbers, then reset the style back to Normal
# (default, unstyled) to match the original (unstyled) B2:B4 cells.
$ws.Range("B2:B4").NumberFormat = "@"

# Row 2 becomes Male / 83.59%
$ws.Range("A2").Value = "Male"
$ws.Range("B2").Value = "83.59%"

# Row 3 becomes Female / 14.49%
$ws.Range("A3").Value = "Female"
$ws.Range("B3").Value = "14.49%"

# Row 4 keeps "Other / Non-Disclosed" label, but percentage updates
$ws.Range("A4").Value = "Other / Non-Disclosed"
$ws.Range("B4").Value = "1.92%"

$ws.Range("B2:B4").Style = "Normal"
